# Updated the existed POM Framework for ADDA and IPMS

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# ---- Sheet 1: rename loginTest -> loginADDA ----
$ws1.Name = "loginADDA"
$ws1.Range("A1:B2").Select()
$ws1.PageSetup.Orientation = 1

# ---- Create invalidLoginTest (inserted right after loginADDA first, so it gets the
#      lower sheetId and ends up as the 3rd tab once loginIPMST is inserted between them) ----
$wsInvalid = $wb.Worksheets.Add($null, $ws1)
$wsInvalid.Name = "invalidLoginTest"
$wsInvalid.Range("A1").Value = "fullName"
$wsInvalid.Range("B1").Value = "password"
$wsInvalid.Range("A2").Value = "archadda_admin"
$wsInvalid.Range("B2").Value = "admin"
$wsInvalid.Columns.Item(1).ColumnWidth = 17
$wsInvalid.Range("D7").Select()

# ---- Create loginIPMST (inserted after loginADDA, pushing invalidLoginTest to 3rd place).
#      This is the last sheet touched/selected, so it ends up the active tab. ----
$wsIpmst = $wb.Worksheets.Add($null, $ws1)
$wsIpmst.Name = "loginIPMST"
$wsIpmst.Range("A1").Value = "fullName"
$wsIpmst.Range("B1").Value = "password"
$wsIpmst.Range("A2").Value = "admin"
$wsIpmst.Range("B2").Value = "admin"
$wsIpmst.Columns.Item(1).ColumnWidth = 17.67
$wsIpmst.Range("B9").Select()

$wb.Save()
